# Apply the budget-estimate bug fix: correct the "booked" hours in B22
# (OneCrossingLandBoundary test fix: 8+34 -> 8+34+34+16), which cascades
# into the dependent "remaining" (B23) and percentage (C23) formulas.
# Also update the active cell selection on the "begroting" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("begroting")
$ws.Activate()

# Update the formula in B22 (booked hours)
$ws.Range("B22").Formula = "=8+34+34+16"

# Recalculate so dependent cells (B23, C23) pick up the new values
$excel.Calculate()

# Move/restore the active cell selection to D20
$ws.Range("D20").Select()
